$d = $word.ActiveDocument

# --- Step 1: remove the trailing "___" paragraph (paragraph 5) ---
$d.Paragraphs(5).Range.Delete()

# --- Step 2: remove the "**_NOTE:_**..." paragraph (paragraph 4); its
#     content is being rewritten into paragraph 2 below ---
$d.Paragraphs(4).Range.Delete()

# --- Step 3: paragraph 3 keeps only its bookmark; drop the text run ---
$p3 = $d.Paragraphs(3)
$p3.Range.Find.Execute("Your blog post should include everything from how you identified what tables contained the information you need, to how you retrieved it using SQL (and any challenges you ran into while doing so), as well as your methodology and results for your hypothesis tests. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# --- Step 4: rewrite paragraph 2 ("### Blog Post Must-Haves") into the
#     multi-run NOTE paragraph ---
$p2xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:i/><w:iCs/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>NOTE:</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> This blog post is your way of showcasing the work you''ve done on this project--chances are it will soon be read by a recruiter or hiring manager! Take the time to make sure that you craft your story well, and clearly explain your process and findings in a way that clearly shows both your technical expertise </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:i/><w:iCs/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t>and</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> your ability to communicate your results!</w:t></w:r></w:p>'
$d.Paragraphs(2).Range.InsertXML($p2xml)

# --- Step 5: rewrite paragraph 1 ("Blog Post") into the long intro text ---
$p1xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>Your blog post should include everything from how you identified what tables contained the information you need, to how you retrieved it using SQL (and any challenges you ran into while doing so), as well as your methodology and results for your hypothesis tests.</w:t></w:r></w:p>'
$d.Paragraphs(1).Range.InsertXML($p1xml)
